# 2017-02-13 snapshot - chunk 30
# Updates the STEO Fig30 "heating degree days" sheet from the January 2017
# edition to the February 2017 edition: refreshed title/source captions and
# refreshed monthly heating degree-day figures for each of the five series
# (2014/15, 2015/16, 2016/17, 2017/18, 2006-2016 Avg).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title and source captions -----------------------------------------
$ws.Range("A2").Value  = "Short-Term Energy Outlook, February 2017"
$ws.Range("A34").Value = "Source: Short-Term Energy Outlook, February 2017."

# --- Monthly heating degree-day data (B27:F32) --------------------------
# Columns: B=2014/15  C=2015/16  D=2016/17  E=2017/18  F=2006-2016 Avg
# Row 27 = October
$ws.Range("B27").Value = 220.44996631999999
$ws.Range("C27").Value = 226.95825633000001
$ws.Range("D27").Value = 197.05462768999999
$ws.Range("E27").Value = 243.15411477999999
$ws.Range("F27").Value = 256.94560739000002

# Row 28 = November
$ws.Range("B28").Value = 613.95444473999999
$ws.Range("C28").Value = 445.21883613
$ws.Range("D28").Value = 417.05579720999998
$ws.Range("E28").Value = 485.11934280999998
$ws.Range("F28").Value = 514.82145361000005

# Row 29 = December
$ws.Range("B29").Value = 705.22987746000001
$ws.Range("C29").Value = 581.17146408999997
$ws.Range("D29").Value = 782.66832546000001
$ws.Range("E29").Value = 767.01030118000006
$ws.Range("F29").Value = 762.39734824000004

# Row 30 = January
$ws.Range("B30").Value = 889.91030916
$ws.Range("C30").Value = 870.11669567000001
$ws.Range("D30").Value = 744.49828066999999
$ws.Range("E30").Value = 844.34603315000004
$ws.Range("F30").Value = 887.55268064999996

# Row 31 = February
$ws.Range("B31").Value = 866.62847237999995
$ws.Range("C31").Value = 627.91217326000003
$ws.Range("D31").Value = 685.60789922000004
$ws.Range("E31").Value = 682.85366234000003
$ws.Range("F31").Value = 746.69470000000001

# Row 32 = March
$ws.Range("B32").Value = 583.53791102000002
$ws.Range("C32").Value = 449.17779856999999
$ws.Range("D32").Value = 551.82551949000003
$ws.Range("E32").Value = 556.64282671000001
$ws.Range("F32").Value = 557.49990000000003

# Row 33 (Total Winter) holds =SUM(...) formulas over B27:F32, so Excel
# recalculates them automatically from the new monthly figures above.
